$wb = $excel.ActiveWorkbook

# --- 1. "Test Cases" overview sheet: flip which test case is active (Runmode) ---
$wsTC = $wb.Worksheets.Item("Test Cases")
$wsTC.Range("C2").Value = "No"   # TC_IncidentReport -> No
$wsTC.Range("C6").Value = "Yes"  # TC_IncidentReportPatientComplaintDetails -> Yes

# --- 2. "TC_InciRepPatientComplaint" sheet: add a "Patient Name" column and fill in data ---
$ws7 = $wb.Worksheets.Item("TC_InciRepPatientComplaint")

# Insert a new column before column D (shifts Complaint DateTime .. Incident IR code right by one)
$ws7.Range("D1").EntireColumn.Insert()

# New column header + value
$ws7.Range("D1").Value = "Patient Name"
$ws7.Range("D2").Value = "Mr. Fashish"

# Update existing values in row 2
$ws7.Range("B2").Value = "Billing"        # Location of Incident: KIMS2 -> Billing
$ws7.Range("J2").Value = "Clinical"       # Department Involved (was I, now J): Client -> Clinical
$ws7.Range("L2").Value = "Eric M Doc"     # Action Taken By (was K, now L): admin -> Eric M Doc
$ws7.Range("N2").Value = "Eric M Doc"     # Witnessed By (was M, now N): admin -> Eric M Doc
$ws7.Range("O2").Value = "Raghu M Doc"    # Notified To (was N, now O): admin -> Raghu M Doc
$ws7.Range("P2").Value = "345 /2022 CONFIG CODE"  # Notified DateTime (was O, now P)

# --- 3. View state: make the patient-complaint sheet the active tab/selection ---
$ws7.Activate()
$ws7.Range("Q2").Select()
